$d = $word.ActiveDocument

# Locate the paragraph that ends with the sock-counting explanation; the new
# "Predicting Fingers" content should be inserted right after it.
$anchorText = "With 8 socks you will have the best chance at the least amount of socks grabbed to get a pair of each. "

$rng = $d.Content
$found = $rng.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate anchor paragraph text"
}

# Build a zero-length insertion point immediately after the anchor paragraph's
# text (but before its paragraph mark) so the new paragraphs land between the
# anchor paragraph and the one that currently follows it.
$insertAt = $rng.End
$insertionPoint = $d.Range($insertAt, $insertAt)

# Two new "List Paragraph"-styled paragraphs: a blank spacer paragraph and a
# paragraph carrying the new "Predicting Fingers:" heading text.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr></w:p>' +
       '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr><w:r><w:t>Predicting Fingers:</w:t></w:r></w:p>' +
       '</w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint.InsertXML($xml) | Out-Null
